# Update "想去人数" (column F) counts on both the "展览" and "全部类型"
# sheets. Rows 2,3,6,7,8,9,10 each get incremented to their new values.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 133
    3  = 1683
    6  = 462
    7  = 153
    8  = 75
    9  = 602
    10 = 408
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
